$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.67 = 26700.0 pesos`n✅ 26700.0 pesos = 6.68 = 960.13 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update N10/O10 and N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 150
$ws2.Range("O10").Value = 4005
$ws2.Range("N12").Value = 3997.5
$ws2.Range("O12").Value = 143.75
